$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_data2")
$ws.Range("C4").Value = "This is test1 table"
